$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ISO9126 note, reused across E11/E12/E13 (becomes a single new shared
# string, matching the dedup that happened in the real edit).
$note = [char]0x2022 + " Utlizes metrics found and based on the ISO9126 guidelines (similar in respect to other models, but different in implementation)"
$compendiumNote = "Note: comprensive metrics for the ISO9126 can be found at http://www.arisa.se/compendium/node6.html"

# Row 11 (Metrics column) was blank before - fill it in and wrap the text so
# it matches the formatting used elsewhere in that column.
$ws.Range("E11").WrapText = $true
$ws.Range("E11").Value = $note

# Rows 12 & 13 (Metrics column) were also blank - same note, same wrap fix.
$ws.Range("E12").WrapText = $true
$ws.Range("E12").Value = $note

$ws.Range("E13").WrapText = $true
$ws.Range("E13").Value = $note

# New row 16: a quick note about where to find the ISO9126 metrics.
$ws.Range("A16").Value = $compendiumNote

# Move the active selection like the author left it.
$ws.Range("B20").Select()
